# Updated symbol list on Mon Jan  2 16:52:57 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) quotes for
# the crypto tracker sheet. Both columns hold plain text (e.g. "246.88",
# "0.91%") rather than numbers, so each new value is written with a
# leading apostrophe to force Excel to keep it as literal text instead of
# re-interpreting it as a number/percentage, and the cell style is reset
# back to "Normal" right afterwards so the quote-prefix flag doesn't leave
# a stray formatting change behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2'  = '246.89'
    'E2'  = '0.85%'
    'D3'  = '29.42'
    'E3'  = '8.48%'
    'D4'  = '5.168'
    'E4'  = '1.84%'
    'D5'  = '0.05711'
    'E5'  = '0.42%'
    'D6'  = '6.598'
    'E6'  = '1.90%'
    'D7'  = '0.8585'
    'D8'  = '0.8716'
    'E8'  = '3.92%'
    'D9'  = '0.1366'
    'E9'  = '2.86%'
    'D10' = '0.07084'
    'E10' = '2.54%'
    'D11' = '0.02870'
    'E11' = '0.75%'
    'D12' = '0.09382'
    'E12' = '-0.17%'
    'D13' = '0.001520'
    'E13' = '-0.69%'
    'D14' = '0.04172'
    'E14' = '1.83%'
    'D15' = '0.0006017'
    'E15' = '0.32%'
    'D16' = '0.006154'
    'E16' = '0.51%'
    'E17' = '3,766.25%'
    'E18' = '-0.70%'
    'D19' = '3.061'
    'E19' = '1.99%'
    'D20' = '2.190'
    'E20' = '-5.19%'
    'E21' = '0.71%'
    'D22' = '0.03319'
    'E22' = '4.42%'
    'E23' = '0.30%'
    'D24' = '3.475'
    'E24' = '-2.73%'
    'E25' = '0.41%'
    'D26' = '0.005040'
    'E26' = '27.55%'
    'E27' = '0.15%'
    'D28' = '0.0001210'
    'E28' = '23.51%'
    'D40' = '0.03755'
    'E40' = '1.73%'
    'D41' = '0.005799'
    'E41' = '-0.65%'
    'E42' = '1.92%'
    'D43' = '0.002540'
    'E43' = '8.45%'
    'E44' = '9.33%'
    'D45' = '0.00005166'
    'E45' = '-0.54%'
    'D47' = '0.06997'
    'E47' = '-31.05%'
    'D48' = '0.002568'
    'E48' = '-0.87%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}

Write-Output "Updated $($updates.Count) cells"
